# TimeSheet.xlsx update — log recent Talos/crawler work on the "Ningge" sheet,
# widen a couple of columns to fit the longer notes, and leave the UI focused
# on the newly-edited sheet/cell (matches the author's own Excel session state).

$wb = $excel.ActiveWorkbook

# --- "Ningge" sheet: fill in the three new timesheet rows -------------------
$wsNingge = $wb.Worksheets.Item("Ningge")

# Row 2 — 2015-10-22, 4.5 hours
$wsNingge.Cells.Item(2, 1).Value = 42299
$wsNingge.Cells.Item(2, 2).Value = 4.5

# Row 3 — 2015-10-23, 2 hours
$wsNingge.Cells.Item(3, 1).Value = 42300
$wsNingge.Cells.Item(3, 2).Value = 2

# Row 4 — 2015-10-25, 4 hours
$wsNingge.Cells.Item(4, 1).Value = 42302
$wsNingge.Cells.Item(4, 2).Value = 4

# Write the "What I worked on" notes in the same order the author typed them
# so the shared-string table comes out in the same sequence as the source file.
$wsNingge.Cells.Item(3, 3).Value = "Change Talos bullet to particle"
$wsNingge.Cells.Item(4, 3).Value = "Add crawler enemy + Crawler animations + crawler bullet animations + crawler shoot bullet towards Talos"
$wsNingge.Cells.Item(2, 3).Value = "Talos animations + Talos controls + Talos bullet animations + Talos shoot bullet to the position where the left button of the mouse being pressed"

# Widen the columns now that there's real content (column C holds long notes).
$wsNingge.Columns.Item(1).ColumnWidth = 12
$wsNingge.Columns.Item(2).ColumnWidth = 25.5
$wsNingge.Columns.Item(3).ColumnWidth = 151.33333333333334

# --- "Team Meetings" sheet: column B needs more room too --------------------
$wsTeam = $wb.Worksheets.Item("Team Meetings")
$wsTeam.Columns.Item(2).ColumnWidth = 24.666666666666668

# --- Leave the selections where the author left them -------------------------
$wsTeam.Range("C20").Select() | Out-Null

$wsProvencher = $wb.Worksheets.Item("Provencher")
$wsProvencher.Range("C17").Select() | Out-Null

# Ningge is the sheet left on-screen/active at the end of the session.
$wsNingge.Activate() | Out-Null
$wsNingge.Range("B8").Select() | Out-Null
